# Applies scheduled market-price refresh values to the Leve profit tables
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*) across all job sheets.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1787.25
$ws.Range("I18").Value = 1854.7
$ws.Range("J18").Value = 1450
$ws.Range("K18").Value = 1854.7
$ws.Range("L18").Value = 1450
$ws.Range("M18").Value = -1570.7
$ws.Range("N18").Value = -2018
$ws.Range("H40").Value = 2624.875
$ws.Range("I40").Value = 1999.8
$ws.Range("J40").Value = 3666.6667
$ws.Range("K40").Value = 1999.8
$ws.Range("L40").Value = 3666.6667
$ws.Range("M40").Value = -1824.8
$ws.Range("N40").Value = -4016.6667
$ws.Range("H74").Value = 2989.8572
$ws.Range("I74").Value = 2738.1667
$ws.Range("K74").Value = 2738.1667
$ws.Range("M74").Value = -1802.1667
$ws.Range("H77").Value = 2989.8572
$ws.Range("I77").Value = 2738.1667
$ws.Range("K77").Value = 13690.8335
$ws.Range("M77").Value = -9010.833500000001
$ws.Range("H92").Value = 970.4286
$ws.Range("I92").Value = 1046.2222
$ws.Range("J92").Value = 834
$ws.Range("K92").Value = 1046.2222
$ws.Range("L92").Value = 834
$ws.Range("M92").Value = 201.7778000000001
$ws.Range("N92").Value = -3330
$ws.Range("H97").Value = 6128
$ws.Range("J97").Value = 6128
$ws.Range("L97").Value = 18384
$ws.Range("N97").Value = -19376
$ws.Range("H100").Value = 1387
$ws.Range("I100").Value = 801.5
$ws.Range("K100").Value = 801.5
$ws.Range("M100").Value = -260.5
$ws.Range("H101").Value = 410.5
$ws.Range("I101").Value = 410.5
$ws.Range("K101").Value = 1231.5
$ws.Range("M101").Value = 390.5
$ws.Range("H106").Value = 25346.1
$ws.Range("J106").Value = 8000
$ws.Range("L106").Value = 8000
$ws.Range("N106").Value = -9262
$ws.Range("H107").Value = 400.18182
$ws.Range("I107").Value = 416.375
$ws.Range("K107").Value = 416.375
$ws.Range("M107").Value = 1503.625

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 803.5714
$ws.Range("I2").Value = 695.8333
$ws.Range("J2").Value = 1450
$ws.Range("K2").Value = 695.8333
$ws.Range("L2").Value = 1450
$ws.Range("M2").Value = -582.8333
$ws.Range("N2").Value = -1676
$ws.Range("H32").Value = 15038.3
$ws.Range("I32").Value = 6067.1143
$ws.Range("J32").Value = 27597.96
$ws.Range("K32").Value = 6067.1143
$ws.Range("L32").Value = 27597.96
$ws.Range("M32").Value = -5780.1143
$ws.Range("N32").Value = -28171.96
$ws.Range("H45").Value = 5161.2
$ws.Range("I45").Value = 1827
$ws.Range("K45").Value = 1827
$ws.Range("M45").Value = -1450
$ws.Range("H63").Value = 5118.6
$ws.Range("I63").Value = 3531
$ws.Range("J63").Value = 7500
$ws.Range("K63").Value = 3531
$ws.Range("L63").Value = 7500
$ws.Range("M63").Value = -2845
$ws.Range("N63").Value = -8872
$ws.Range("H66").Value = 5118.6
$ws.Range("I66").Value = 3531
$ws.Range("J66").Value = 7500
$ws.Range("K66").Value = 17655
$ws.Range("L66").Value = 37500
$ws.Range("M66").Value = -14223
$ws.Range("N66").Value = -44364
$ws.Range("H116").Value = 803.5714
$ws.Range("I116").Value = 695.8333
$ws.Range("J116").Value = 1450
$ws.Range("K116").Value = 695.8333
$ws.Range("L116").Value = 1450
$ws.Range("M116").Value = 1598.1667
$ws.Range("N116").Value = -6038
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = $null
$ws.Range("H132").Value = 1290.3182
$ws.Range("I132").Value = 1336.3414
$ws.Range("J132").Value = 661.3333
$ws.Range("K132").Value = 4009.0242
$ws.Range("L132").Value = 1983.9999
$ws.Range("M132").Value = -1479.0242
$ws.Range("N132").Value = -7043.9999

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 803.5714
$ws.Range("I3").Value = 695.8333
$ws.Range("J3").Value = 1450
$ws.Range("K3").Value = 695.8333
$ws.Range("L3").Value = 1450
$ws.Range("M3").Value = -581.8333
$ws.Range("N3").Value = -1678
$ws.Range("H5").Value = 1579.8
$ws.Range("I5").Value = 1724.75
$ws.Range("K5").Value = 1724.75
$ws.Range("M5").Value = -1611.75
$ws.Range("H22").Value = 610.7857
$ws.Range("I22").Value = 623.1539
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 623.1539
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -450.1539
$ws.Range("N22").Value = -796
$ws.Range("H82").Value = 25772.8
$ws.Range("I82").Value = 7216.25
$ws.Range("K82").Value = 7216.25
$ws.Range("M82").Value = -6833.25
$ws.Range("H85").Value = 25772.8
$ws.Range("I85").Value = 7216.25
$ws.Range("K85").Value = 7216.25
$ws.Range("M85").Value = -5890.25
$ws.Range("H94").Value = 1199
$ws.Range("I94").Value = 998.75
$ws.Range("K94").Value = 998.75
$ws.Range("M94").Value = -547.75
$ws.Range("H99").Value = 1250
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1250
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 1250
$ws.Range("M99").Value = $null
$ws.Range("N99").Value = -4246
$ws.Range("H105").Value = 3836.516
$ws.Range("I105").Value = 3151.7
$ws.Range("K105").Value = 3151.7
$ws.Range("M105").Value = -1404.7
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
$ws.Range("H135").Value = 93278
$ws.Range("J135").Value = 93278
$ws.Range("L135").Value = 93278
$ws.Range("N135").Value = -103418

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4203.2915
$ws.Range("I31").Value = 2694.6365
$ws.Range("K31").Value = 2694.6365
$ws.Range("M31").Value = -2399.6365
$ws.Range("H34").Value = 4203.2915
$ws.Range("I34").Value = 2694.6365
$ws.Range("K34").Value = 2694.6365
$ws.Range("M34").Value = -2492.6365
$ws.Range("H105").Value = 3307
$ws.Range("I105").Value = 2277
$ws.Range("J105").Value = 4079.5
$ws.Range("K105").Value = 2277
$ws.Range("L105").Value = 4079.5
$ws.Range("M105").Value = -530
$ws.Range("N105").Value = -7573.5
$ws.Range("H122").Value = 8746
$ws.Range("I122").Value = 8783.571
$ws.Range("J122").Value = 8680.25
$ws.Range("K122").Value = 26350.713
$ws.Range("L122").Value = 26040.75
$ws.Range("M122").Value = -23900.713
$ws.Range("N122").Value = -30940.75
$ws.Range("H132").Value = 3302.75
$ws.Range("I132").Value = 2583.25
$ws.Range("K132").Value = 7749.75
$ws.Range("M132").Value = -5219.75
$ws.Range("H134").Value = 4329
$ws.Range("I134").Value = 3178
$ws.Range("K134").Value = 9534
$ws.Range("M134").Value = -6999

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1918.5
$ws.Range("J68").Value = 1792
$ws.Range("L68").Value = 5376
$ws.Range("N68").Value = -6998
$ws.Range("H71").Value = 1918.5
$ws.Range("J71").Value = 1792
$ws.Range("L71").Value = 16128
$ws.Range("N71").Value = -24240
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = $null
$ws.Range("N130").Value = $null
$ws.Range("H131").Value = 4092.3333
$ws.Range("I131").Value = 3022.111
$ws.Range("K131").Value = 9066.332999999999
$ws.Range("M131").Value = -4026.332999999999
$ws.Range("H140").Value = 3373.5
$ws.Range("I140").Value = 3373.5
$ws.Range("K140").Value = 10120.5
$ws.Range("M140").Value = -4940.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 44996.5
$ws.Range("J100").Value = 44996.5
$ws.Range("L100").Value = 44996.5
$ws.Range("N100").Value = -47160.5

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2041.1428
$ws.Range("I7").Value = 2041.1428
$ws.Range("K7").Value = 2041.1428
$ws.Range("M7").Value = -1929.1428
$ws.Range("H22").Value = 1501
$ws.Range("J22").Value = 1002
$ws.Range("L22").Value = 1002
$ws.Range("N22").Value = -1592
$ws.Range("H27").Value = 1501
$ws.Range("J27").Value = 1002
$ws.Range("L27").Value = 1002
$ws.Range("N27").Value = -1216
$ws.Range("H40").Value = 2326.8572
$ws.Range("I40").Value = 2326.8572
$ws.Range("K40").Value = 2326.8572
$ws.Range("M40").Value = -2190.8572
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = $null
$ws.Range("H122").Value = 950
$ws.Range("J122").Value = 950
$ws.Range("L122").Value = 2850
$ws.Range("N122").Value = -7750
$ws.Range("H126").Value = 2041.1428
$ws.Range("I126").Value = 2041.1428
$ws.Range("K126").Value = 6123.428400000001
$ws.Range("M126").Value = -3653.428400000001
$ws.Range("H127").Value = 86499.336
$ws.Range("J127").Value = 86499.336
$ws.Range("L127").Value = 86499.336
$ws.Range("N127").Value = -96419.336

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 19999
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = $null
$ws.Range("H62").Value = 7624.375
$ws.Range("J62").Value = 7624.375
$ws.Range("L62").Value = 7624.375
$ws.Range("N62").Value = -8872.375
$ws.Range("H65").Value = 7624.375
$ws.Range("J65").Value = 7624.375
$ws.Range("L65").Value = 38121.875
$ws.Range("N65").Value = -44361.875
$ws.Range("H96").Value = 1636.5555
$ws.Range("I96").Value = 1609.8
$ws.Range("J96").Value = 1670
$ws.Range("K96").Value = 1609.8
$ws.Range("L96").Value = 1670
$ws.Range("M96").Value = -236.8
$ws.Range("N96").Value = -4416
$ws.Range("I100").Value = 1129.2
$ws.Range("J100").Value = 1595
$ws.Range("K100").Value = 2258.4
$ws.Range("L100").Value = 3190
$ws.Range("M100").Value = -1717.4
$ws.Range("N100").Value = -4272
$ws.Range("H113").Value = 1505.6364
$ws.Range("I113").Value = 1447.9
$ws.Range("J113").Value = 1553.75
$ws.Range("K113").Value = 4343.700000000001
$ws.Range("L113").Value = 4661.25
$ws.Range("M113").Value = -2173.700000000001
$ws.Range("N113").Value = -9001.25
$ws.Range("H122").Value = 1312.2142
$ws.Range("I122").Value = 1312.2142
$ws.Range("K122").Value = 3936.6426
$ws.Range("M122").Value = -1486.6426
$ws.Range("H126").Value = 57555
$ws.Range("I126").Value = 63686.875
$ws.Range("K126").Value = 191060.625
$ws.Range("M126").Value = -188590.625
$ws.Range("H132").Value = 1361.05
$ws.Range("I132").Value = 1311.8
$ws.Range("K132").Value = 3935.4
$ws.Range("M132").Value = -1405.4

Write-Host "Applied market price updates"
